# Replace the stack-trace text of the setNumbering error run.
# The original run's <w:t> contains a long, literal multi-line Java stack
# trace (using real \n / \t characters, not <w:br/> elements). We locate
# the run by its distinctive start/end markers (scoped to the second
# paragraph, which holds the whole template body) and then overwrite the
# full span with the updated stack trace text (AQL 7.x -> 8.x regeneration),
# including the trailing tab character that now precedes the closing tag.

$d = $word.ActiveDocument
$p2 = $d.Paragraphs(2)

$rStart = $d.Content
$rStart.Start = $p2.Range.Start
$rStart.End = $p2.Range.End
$foundStart = $rStart.Find.Execute("<---setNumbering", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rEnd = $d.Content
$rEnd.Start = $p2.Range.Start
$rEnd.End = $p2.Range.End
$foundEnd = $rEnd.Find.Execute("RemoteTestRunner.main(RemoteTestRunner.java:210)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundStart -or -not $foundEnd) {
    Write-Host "ERROR: could not locate stack trace markers (start=$foundStart end=$foundEnd)"
} else {
    $newText = '<---setNumbering(org.obeonetwork.m2doc.element.MParagraph,java.lang.Integer,java.lang.Integer) with arguments [org.obeonetwork.m2doc.element.impl.MParagraphImpl@60cc20e1, 9999, 3] failed:
	java.lang.IllegalArgumentException: no numbering with ID 9999
		at org.obeonetwork.m2doc.services.PaginationServices.setNumbering(PaginationServices.java:414)
		at java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)
		at java.base/java.lang.reflect.Method.invoke(Method.java:580)
		at org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)
		at org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallOrApply(EvaluationServices.java:273)
		at org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:173)
		at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)
		at org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)
		at org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:674)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
		at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
		at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)
		at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)
		at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)
		at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)
		at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)
		at java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)
		at java.base/java.lang.reflect.Method.invoke(Method.java:580)
		at org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59)
		at org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)
		at org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)
		at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
		at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
		at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)
		at org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100)
		at org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)
		at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)
		at org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)
		at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)
		at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)
		at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)
		at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)
		at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)
		at org.junit.runners.ParentRunner.run(ParentRunner.java:413)
		at org.junit.runners.Suite.runChild(Suite.java:128)
		at org.junit.runners.Suite.runChild(Suite.java:27)
		at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)
		at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)
		at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)
		at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)
		at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)
		at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
		at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
		at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)
		at org.junit.runners.ParentRunner.run(ParentRunner.java:413)
		at org.junit.runners.Suite.runChild(Suite.java:128)
		at org.junit.runners.Suite.runChild(Suite.java:27)
		at org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)
		at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)
		at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)
		at org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)
		at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)
		at org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)
		at org.junit.runners.ParentRunner.run(ParentRunner.java:413)
		at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)
		at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)
		at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)
	'

    # +1 to also include the single trailing newline character that sits
    # right after the marker text and right before the <w:br/> break.
    $target = $d.Range($rStart.Start, $rEnd.End + 1)
    $target.Text = $newText

    Write-Host "Replaced stack trace text: $($newText.Length) chars"
}
